$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card18")

# Row 17 used to be the "open" trailing row with blank placeholder cells in
# columns B-K and N. Now that a new event (row 18) is being appended, row 17
# becomes a "closed" row and its placeholder cells are finalized to the
# literal text "nan" (matching every other fully-populated row in the table).
$ws.Range("B17").Value = "nan"
$ws.Range("C17").Value = "nan"
$ws.Range("D17").Value = "nan"
$ws.Range("E17").Value = "nan"
$ws.Range("F17").Value = "nan"
$ws.Range("G17").Value = "nan"
$ws.Range("H17").Value = "nan"
$ws.Range("I17").Value = "nan"
$ws.Range("J17").Value = "nan"
$ws.Range("K17").Value = "nan"
$ws.Range("N17").Value = "nan"

# Row 18 is the new event row: "card" number plus the new service entry.
# Columns B-K and N stay blank (present-but-empty cells), mirroring the
# previous "open" row pattern, while L/M/O/P carry the new data.
$ws.Range("A18").Value = "'18"

$ws.Range("B18:K18").Style = "Normal"
$ws.Range("N18").Style = "Normal"

$ws.Range("L18").Value = "14\8\2025"
$ws.Range("M18").Value = "9736 h   775 t"
$ws.Range("O18").Value = "تم تغيير زيت الجيربوكس"
$ws.Range("P18").Value = "تيم العمل"
